$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prefix the hotel id (00103) onto id_servicios values in column A (rows 2-11)
$values = @(1000310101, 1000310102, 1000310103, 1000310104, 1000310105, 1000310106, 1000310107, 1000310108, 1000310109, 1000310110)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Widen column A so the longer ids fit (best-fit sized)
$ws.Columns.Item(1).ColumnWidth = 10.330729166666666

# Selection now spans the updated id_servicios column, active cell on the first data row
$ws.Range("A2:A11").Select()
